$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the two case descriptions to reflect the new deterministic runs
$ws.Range("C14").Value = "Age and length comps, deterministic"
$ws.Range("C15").Value = "Length and calcomps, deterministic"

# Move the active selection to C4, matching the edited view state
$ws.Range("C4").Select()
